# Switch to using only the first row for Date and Opponent_Team.
# The C/D columns (Date, Opponent_Team) were being repeated on every
# row of each results sheet; only the first data row (row 2) needs to
# keep the value now, so clear the redundant copies in rows 3-7.

$wb = $excel.ActiveWorkbook

$sheetNames = @("results-1", "results-2", "results-3")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("C3:D7").ClearContents()
}

# Restore the selection / active-cell state recorded for each sheet,
# ending with results-3 active (matches tabSelected on that sheet).
$ws1 = $wb.Worksheets.Item("results-1")
$ws1.Activate()
$ws1.Range("E16").Select()

$ws2 = $wb.Worksheets.Item("results-2")
$ws2.Activate()
$ws2.Range("F17").Select()

$ws3 = $wb.Worksheets.Item("results-3")
$ws3.Activate()
$ws3.Range("F20").Select()
